$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("C1").Value = "event_name"
$ws.Range("D1").Value = "result_month"
$ws.Range("E1").Value = "event_result"

# --- Column C: event_name (same value repeated for every data row) ---
$eventName = "Разметка данных на CVAT для начинающих"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = $eventName
}

# --- Column D: result_month, left blank for every data row ---
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 4).Value = ""
}

# --- Column E: event_result values (previously column F) ---
$ws.Cells.Item(2, 5).Value = "Сертификат участника"
$ws.Cells.Item(3, 5).Value = "Сертификат участника"
$ws.Cells.Item(4, 5).Value = "Сертификат участника"
$ws.Cells.Item(5, 5).Value = " "
$ws.Cells.Item(6, 5).Value = " "
$ws.Cells.Item(7, 5).Value = " "
$ws.Cells.Item(8, 5).Value = " "

# --- Remove the now-obsolete column F entirely ---
$ws.Columns.Item(6).Delete()
